$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) stays text-typed, matching the source
# inlineStr cells exactly (avoids Excel auto-converting numeric-looking
# strings like "1.00" / "0.400" into numbers and dropping trailing zeros).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.344.37"
$ws.Range("E2").Value = "  +0.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.606.25"
$ws.Range("E3").Value = "  -0.11%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.13"
$ws.Range("E5").Value = "  -1.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.32"
$ws.Range("E6").Value = "  +0.67%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  -1.61%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.601.10"
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.183"
$ws.Range("E10").Value = "  +3.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.668"
$ws.Range("E11").Value = "  +1.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.27"
$ws.Range("E12").Value = "  -3.87%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000313"
$ws.Range("E13").Value = "  +8.70%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.74"
$ws.Range("E14").Value = "  -1.15%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.171.62"
$ws.Range("E15").Value = "  -0.21%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.03"
$ws.Range("E16").Value = "  +2.51%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.592.08"
$ws.Range("E17").Value = "  -0.30%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.213.54"
$ws.Range("E18").Value = "  +0.16%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.77"
$ws.Range("E19").Value = "  +1.74%  "

# Row 20
$ws.Range("E20").Value = "  +0.19%  "

# Row 21
$ws.Range("E21").Value = "  -0.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "485.13"
$ws.Range("E22").Value = "  -1.38%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.23"
$ws.Range("E23").Value = "  +11.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.06"
$ws.Range("E24").Value = "  -6.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.41"
$ws.Range("E25").Value = "  -1.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.87"
$ws.Range("E26").Value = "  +5.76%  "

# Row 27
$ws.Range("E27").Value = "  -3.33%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.23"
$ws.Range("E28").Value = "  +1.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.49"
$ws.Range("E29").Value = "  +1.11%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.42"
$ws.Range("E30").Value = "  -0.98%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.67"
$ws.Range("E31").Value = "  +1.34%  "

# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.32"
$ws.Range("E32").Value = "  +0.00%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.120"
$ws.Range("E33").Value = "  +1.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.37"
$ws.Range("E34").Value = "  +1.69%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "583.81"
$ws.Range("E35").Value = "  -6.12%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.22"
$ws.Range("E36").Value = "  +2.56%  "

# Row 37
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0816"
$ws.Range("E37").Value = "  -0.66%  "

# Row 38
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.17%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.400"
$ws.Range("E39").Value = "  -0.95%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.34"
$ws.Range("E40").Value = "  +22.73%  "

# Row 41
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("E41").Value = "  +6.17%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.46"
$ws.Range("E42").Value = "  -3.13%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.137"
$ws.Range("E43").Value = "  -6.55%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.231.08"
$ws.Range("E44").Value = "  -3.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.10"
$ws.Range("E45").Value = "  -3.37%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0452"
$ws.Range("E46").Value = "  +0.23%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.59"
$ws.Range("E47").Value = "  +5.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.37"
$ws.Range("E48").Value = "  +2.68%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.139"
$ws.Range("E49").Value = "  +0.80%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.13%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.15"
$ws.Range("E51").Value = "  -4.03%  "
